$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Properties_Predictive")

# New column H, mirroring the existing F (PredictiveDatasourceName /
# TeradataDS) and G (matrixColorRelFormula / LIB_ISP.updateMatrixColorRelation)
# columns: H1 is a blank title-bar cell, H2 is the bold header
# "DataFieldMatrixFormula", H3 is the data value
# "LIB_ISP.getDataFieldListForMatrix".

# H1 - blank cell of the red title bar (same look as F1/G1).
$ws.Range("H1").Font.Name = "Trebuchet MS"
$ws.Range("H1").Font.Size = 10
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").Font.Color = 16777215
$ws.Range("H1").Interior.Color = 2368940
$ws.Range("H1").WrapText = $true

# H2 - bold white-on-red column header (same look as A2:E2).
$ws.Range("H2").Value = "DataFieldMatrixFormula"
$ws.Range("H2").Font.Name = "Trebuchet MS"
$ws.Range("H2").Font.Size = 10
$ws.Range("H2").Font.Bold = $true
$ws.Range("H2").Font.Color = 16777215
$ws.Range("H2").Interior.Color = 2368940
$ws.Range("H2").WrapText = $true

# H3 - plain Trebuchet MS 10pt body value.
$ws.Range("H3").Value = "LIB_ISP.getDataFieldListForMatrix"
$ws.Range("H3").Font.Name = "Trebuchet MS"
$ws.Range("H3").Font.Size = 10

# Column width / row height to match the new content.
$ws.Columns.Item(8).ColumnWidth = 30.42578125
$ws.Rows.Item(3).RowHeight = 15.75

# Selection moves onto the newly added column.
[void]$ws.Range("H1:H3").Select()
